$d = $word.ActiveDocument

# --- Paragraph 1: "Statement for "About Us" section of website." ---
# Wrap the whole sentence in gramStart/gramEnd proofErr marks (Word's grammar
# checker flagged it), keeping the existing two runs and paragraph
# properties intact.
$p1 = $d.Paragraphs(1)
if ($p1.Range.Text -notlike "Statement for*section of website.*") {
    throw "Paragraph 1 text did not match expectations: $($p1.Range.Text)"
}
$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>Statement for &#8220;About Us&#8221; section of website</w:t></w:r><w:r><w:t>.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@
$null = $p1.Range.InsertXML($xml1)

# --- Paragraph "Agreed to restaurant ... message board." ---
# Split into two sentences/runs and wrap the first sentence with
# gramStart/gramEnd proofErr marks.
$pAgreed = $d.Paragraphs(5)
if ($pAgreed.Range.Text -notlike "Agreed to restaurant*message board.*") {
    throw "Paragraph 5 text did not match expectations: $($pAgreed.Range.Text)"
}
$xmlAgreed = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>Agreed to restaurant with delivery service available to be our website.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> This was decided on our group&#8217;s message board.</w:t></w:r></w:p>
'@
$null = $pAgreed.Range.InsertXML($xmlAgreed)

# --- Paragraph "Meeting: 5/31/14 2:00pm PDT, duration: 1 hour. All participated. ..." ---
# Split the first run into three runs, inserting the GMT time, and leave the
# remaining (unchanged) runs of the paragraph as-is.
$pMeeting = $d.Paragraphs(6)
if ($pMeeting.Range.Text -notlike "Meeting: 5/31/14*too complicated.*") {
    throw "Paragraph 6 text did not match expectations: $($pMeeting.Range.Text)"
}
$xmlMeeting = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Meeting: 5/31/14 2:00pm PDT</w:t></w:r><w:r><w:t xml:space="preserve"> (21:00 GMT)</w:t></w:r><w:r><w:t>, duration: 1 hour</w:t></w:r><w:r><w:t xml:space="preserve">. All participated. </w:t></w:r><w:r><w:t xml:space="preserve">Decided on sandwich shop since allowing too many options would make the menu and website </w:t></w:r><w:r><w:t>overall</w:t></w:r><w:r><w:t xml:space="preserve"> too complicated.</w:t></w:r></w:p>
'@
$null = $pMeeting.Range.InsertXML($xmlMeeting)

Write-Output $d.Content.Text
